$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column B (rows 2 to 15) with the insulation type name
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 2).Value = "Isover Tapelock Rørskål"
}

# Set column B width to fit content (bestFit); 21.1666... (= 22 - 5/6) is the
# ColumnWidth value that round-trips to a stored OOXML width of exactly 22.
$ws.Columns.Item(2).ColumnWidth = 21.166666666666668

# Update the active selection to C18
$ws.Range("C18").Select()
